$wb = $excel.ActiveWorkbook

$oldVersion = "mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)"
$newVersion = "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on February 03 2026 17.29.55 EST)"

# --- "About" sheet updates ---
$about = $wb.Worksheets.Item("About")

$about.Range("A2").Value = "Version: " + $newVersion

$newCitation = "Recommended Citation:  `"Global Energy Monitor, Coal mine boundaries and methane sources for Vorkutinskaya Coal Mine, Russia, M0868, version '" + $newVersion + "'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

$about.Range("A6").Value = $newCitation

# --- "Boundaries and methane sources" sheet updates ---
$data = $wb.Worksheets.Item("Boundaries and methane sources")

$usedRange = $data.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $data.Cells.Item($r, 19)  # column S = 19
    if ($cell.Text -eq $oldVersion) {
        $cell.Value = $newVersion
    }
}

$wb.Save()
